$wb = $excel.ActiveWorkbook

# --- 1. Rename sheets ---
$wsDaily  = $wb.Worksheets.Item("daily")
$wsMinute = $wb.Worksheets.Item("minute")
$wsSheet1 = $wb.Worksheets.Item("Sheet1")

$wsDaily.Name  = "d"
$wsMinute.Name = "m1"
$wsSheet1.Name = "h1"

# --- 2. Add new sheet "h4" right after "h1" ---
$wsH4 = $wb.Worksheets.Add($null, $wsSheet1)
$wsH4.Name = "h4"

Write-Output "sheets renamed/added"

# --- 3. Populate ticker lists ---
# Register the two brand-new shared strings (USD_CHF, USD_CAD) in that exact
# order first, so they land at shared-string indices 8 and 9 respectively -
# matching the author's original edit order - before laying out the final,
# alphabetically sorted column on sheet "d".
$wsD = $wb.Worksheets.Item("d")
$wsD.Cells.Item(9,1).Value  = "USD_CHF"
$wsD.Cells.Item(8,1).Value  = "USD_CAD"

$sortedTickers = @("Tickers","AUD_USD","EUR_JPY","EUR_USD","GBP_JPY","GBP_USD","NZD_USD","USD_CAD","USD_CHF","USD_JPY")

$wsM1 = $wb.Worksheets.Item("m1")
$wsH1 = $wb.Worksheets.Item("h1")

for ($i = 0; $i -lt $sortedTickers.Length; $i++) {
    $r = $i + 1
    $wsD.Cells.Item($r,1).Value  = $sortedTickers[$i]
    $wsM1.Cells.Item($r,1).Value = $sortedTickers[$i]
    $wsH1.Cells.Item($r,1).Value = $sortedTickers[$i]
}

$wsH4 = $wb.Worksheets.Item("h4")
$wsH4.Cells.Item(1,1).Value = "Tickers"
$wsH4.Cells.Item(2,1).Value = "EUR_USD"
$wsH4.Cells.Item(3,1).Value = "GBP_USD"

Write-Output "data written"

# --- 4. AutoFilter + hidden _FilterDatabase name on sheet "d" ---
$wsD.Range("A1:A10").AutoFilter() | Out-Null
$filterName = $wsD.Names.Add("_xlnm._FilterDatabase", "=d!`$A`$1:`$A`$10")
$filterName.Visible = $false

Write-Output "autofilter set"

# --- 5. Conditional formatting: highlight duplicate tickers ---
# dxfId allocation follows creation order, so rules are added in the same
# order the author must have added them (h4, then h1 x2, then m1, then d)
# to land on dxfId 0,1,2,3,4 respectively.
function Add-DuplicateHighlight($ws, $rangeAddr) {
    $rng = $ws.Range($rangeAddr)
    $fc = $rng.FormatConditions.AddUniqueValues()
    $fc.DupeUnique = 1
    $fc.Font.Color = 393372
    $fc.Interior.Color = 13551615
}

Add-DuplicateHighlight $wsH4 "A4:A10"

Add-DuplicateHighlight $wsH1 "A2:A10"
Add-DuplicateHighlight $wsH1 "A11:A1048576"

Add-DuplicateHighlight $wsM1 "A2:A10"

Add-DuplicateHighlight $wsD "A2:A10"

Write-Output "conditional formatting set"

# --- 6. Final selections / active sheet ---
$wsD.Range("D21").Select() | Out-Null
$wsM1.Range("A1:A10").Select() | Out-Null
$wsH1.Range("A1:A10").Select() | Out-Null
$wsH4.Range("F11").Select() | Out-Null

Write-Output "selections set"
